$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header B1 to mention the new "league" match type
$ws.Range("B1").Value = "Match Type (casual/tournament/league)"

# 2. Row 3 (singles tournament match) was missing the Player DOB columns
#    (E and H). Capture the existing values first, then re-write the row
#    with the two new DOB values inserted, shifting the trailing columns
#    back into their correct (schema-aligned) positions.
$r3_N = $ws.Range("N3").Value2
$r3_O = $ws.Range("O3").Value2
$r3_P = $ws.Range("P3").Value2
$r3_Q = $ws.Range("Q3").Value2
$r3_R = $ws.Range("R3").Value2
$r3_S = $ws.Range("S3").Value2
$r3_T = $ws.Range("T3").Value2
$r3_W = $ws.Range("W3").Value2
$r3_X = $ws.Range("X3").Value2

$ws.Range("N3").Value = ""
$ws.Range("W3").Value = ""

$ws.Range("E3").Value = 31116
$ws.Range("H3").Value = 33076

$ws.Range("O3").Value = $r3_N
$ws.Range("P3").Value = $r3_O
$ws.Range("Q3").Value = $r3_P
$ws.Range("R3").Value = $r3_Q
$ws.Range("S3").Value = $r3_R
$ws.Range("T3").Value = $r3_S
$ws.Range("U3").Value = $r3_T
$ws.Range("X3").Value = $r3_W
$ws.Range("Y3").Value = $r3_X

# 3. Row 4 (close singles match) has the same missing-DOB issue, plus a
#    third game (V/W) that also needs to shift right by one column.
$r4_N = $ws.Range("N4").Value2
$r4_O = $ws.Range("O4").Value2
$r4_P = $ws.Range("P4").Value2
$r4_Q = $ws.Range("Q4").Value2
$r4_R = $ws.Range("R4").Value2
$r4_S = $ws.Range("S4").Value2
$r4_T = $ws.Range("T4").Value2
$r4_U = $ws.Range("U4").Value2
$r4_V = $ws.Range("V4").Value2
$r4_W = $ws.Range("W4").Value2
$r4_X = $ws.Range("X4").Value2

$ws.Range("N4").Value = ""
$ws.Range("V4").Value = ""

$ws.Range("E4").Value = 32452
$ws.Range("H4").Value = 33865

$ws.Range("O4").Value = $r4_N
$ws.Range("P4").Value = $r4_O
$ws.Range("Q4").Value = $r4_P
$ws.Range("R4").Value = $r4_Q
$ws.Range("S4").Value = $r4_R
$ws.Range("T4").Value = $r4_S
$ws.Range("U4").Value = $r4_T
$ws.Range("V4").Value = $r4_U
$ws.Range("W4").Value = $r4_V
$ws.Range("X4").Value = $r4_W
$ws.Range("Y4").Value = $r4_X

# 4. Row 5 is retagged as a "league" match.
$ws.Range("B5").Value = "league"
$ws.Range("Y5").Value = "League doubles match - Team 2 wins in straight sets"
